# Update ACCESS Variable Dictionary 2018-2019.xlsx
# Renames the four "mode_*" variables to "cbt_*" (computer-based testing mode
# indicator instead of online/paper mode) and updates their Values text from
# "2 levels: Online, Paper" to "2 levels: Y, N". The accommodation row's
# Values text is updated the same way, from "2 levels: 1 = Yes, 0 = No" to
# "2 levels: Y, N".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 34-37: mode_listen/read/speak/write -> cbt_listen/read/speak/write
$ws.Range("B34").Value = "cbt_listen"
$ws.Range("B35").Value = "cbt_read"
$ws.Range("B36").Value = "cbt_speak"
$ws.Range("B37").Value = "cbt_write"

# Row 33: accommodation - just the Values (F) column changes.
# Rows 34-37: Values column changes from "2 levels: Online, Paper" to
# "2 levels: Y, N" as well.
$ws.Range("F33").Value = "2 levels: Y, N"
$ws.Range("F34").Value = "2 levels: Y, N"
$ws.Range("F35").Value = "2 levels: Y, N"
$ws.Range("F36").Value = "2 levels: Y, N"
$ws.Range("F37").Value = "2 levels: Y, N"

# Update the saved selection/view state to match (active cell E53, no
# special scroll position).
$ws.Range("E53").Select()
